$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A96").Value = 3234
$ws.Range("B96").Value = "Count the Number of Substrings With Dominant Ones"
$ws.Range("C96").Value = "Math/Loop"

$ws.Range("A97").Value = 1513
$ws.Range("B97").Value = "Number of Substrings With Only 1s"
$ws.Range("C97").Value = "Math/Loop"
$ws.Range("D97").Value = "Triangular Number, sum of sequence of 1s ( 1,2,3,…,n) = n(n+1)/2"

$ws.Range("D97").Select()
